$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 33.65726433333334
$ws.Range("H2").Value = 100.971793
$ws.Range("I2").Value = 0.8115737688004754
$ws.Range("J2").Value = 0.8115737688004754
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.4310329999999999
$ws.Range("N2").Value = 1.293099
$ws.Range("O2").Value = 0.08359060417869307
$ws.Range("P2").Value = 0.08359060417869306
$ws.Range("Q2").Value = 14.50739161738967
$ws.Range("R2").Value = 130.566524556507
$ws.Range("S2").Value = 0.0678399416696107
$ws.Range("T2").Value = 0.06783994166961069
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 33.65726433333334
$ws.Range("H3").Value = 100.971793
$ws.Range("I3").Value = 0.8115737688004754
$ws.Range("J3").Value = 0.8115737688004754
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.734602
$ws.Range("N3").Value = 11.203806
$ws.Range("O3").Value = 0.7242546105447971
$ws.Range("P3").Value = 0.7242546105447971
$ws.Range("Q3").Value = 125.6964866937954
$ws.Range("R3").Value = 1131.268380244158
$ws.Range("S3").Value = 0.5877860438509616
$ws.Range("T3").Value = 0.5877860438509616
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 33.65726433333334
$ws.Range("H4").Value = 100.971793
$ws.Range("I4").Value = 0.8115737688004754
$ws.Range("J4").Value = 0.8115737688004754
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.9908416666666667
$ws.Range("N4").Value = 2.972525
$ws.Range("O4").Value = 0.1921547852765099
$ws.Range("P4").Value = 0.1921547852765099
$ws.Range("Q4").Value = 33.34901988748056
$ws.Range("R4").Value = 300.141178987325
$ws.Range("S4").Value = 0.1559477832799032
$ws.Range("T4").Value = 0.1559477832799032
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.56955
$ws.Range("H5").Value = 1.70865
$ws.Range("I5").Value = 0.01373349406661455
$ws.Range("J5").Value = 0.01373349406661455
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4310329999999999
$ws.Range("N5").Value = 1.293099
$ws.Range("O5").Value = 0.08359060417869307
$ws.Range("P5").Value = 0.08359060417869306
$ws.Range("Q5").Value = 0.24549484515
$ws.Range("R5").Value = 2.209453606349999
$ws.Range("S5").Value = 0.001147991066512806
$ws.Range("T5").Value = 0.001147991066512806
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.56955
$ws.Range("H6").Value = 1.70865
$ws.Range("I6").Value = 0.01373349406661455
$ws.Range("J6").Value = 0.01373349406661455
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.734602
$ws.Range("N6").Value = 11.203806
$ws.Range("O6").Value = 0.7242546105447971
$ws.Range("P6").Value = 0.7242546105447971
$ws.Range("Q6").Value = 2.1270425691
$ws.Range("R6").Value = 19.1433831219
$ws.Range("S6").Value = 0.009946546396635201
$ws.Range("T6").Value = 0.009946546396635201
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.56955
$ws.Range("H7").Value = 1.70865
$ws.Range("I7").Value = 0.01373349406661455
$ws.Range("J7").Value = 0.01373349406661455
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9908416666666667
$ws.Range("N7").Value = 2.972525
$ws.Range("O7").Value = 0.1921547852765099
$ws.Range("P7").Value = 0.1921547852765099
$ws.Range("Q7").Value = 0.56433387125
$ws.Range("R7").Value = 5.07900484125
$ws.Range("S7").Value = 0.002638956603466541
$ws.Range("T7").Value = 0.002638956603466541
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.244787666666667
$ws.Range("H8").Value = 21.734363
$ws.Range("I8").Value = 0.1746927371329101
$ws.Range("J8").Value = 0.1746927371329101
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.4310329999999999
$ws.Range("N8").Value = 1.293099
$ws.Range("O8").Value = 0.08359060417869307
$ws.Range("P8").Value = 0.08359060417869306
$ws.Range("Q8").Value = 3.122742562326333
$ws.Range("R8").Value = 28.104683060937
$ws.Range("S8").Value = 0.01460267144256956
$ws.Range("T8").Value = 0.01460267144256956
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.244787666666667
$ws.Range("H9").Value = 21.734363
$ws.Range("I9").Value = 0.1746927371329101
$ws.Range("J9").Value = 0.1746927371329101
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.734602
$ws.Range("N9").Value = 11.203806
$ws.Range("O9").Value = 0.7242546105447971
$ws.Range("P9").Value = 0.7242546105447971
$ws.Range("Q9").Value = 27.05639850950867
$ws.Range("R9").Value = 243.507586585578
$ws.Range("S9").Value = 0.1265220202972004
$ws.Range("T9").Value = 0.1265220202972004
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.244787666666667
$ws.Range("H10").Value = 21.734363
$ws.Range("I10").Value = 0.1746927371329101
$ws.Range("J10").Value = 0.1746927371329101
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.9908416666666667
$ws.Range("N10").Value = 2.972525
$ws.Range("O10").Value = 0.1921547852765099
$ws.Range("P10").Value = 0.1921547852765099
$ws.Range("Q10").Value = 7.178437486286112
$ws.Range("R10").Value = 64.60593737657501
$ws.Range("S10").Value = 0.03356804539314012
$ws.Range("T10").Value = 0.03356804539314012
